$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Store (B) and Item ID (D) values for rows 2-4
$ws.Range("B2").Value = 5
$ws.Range("D2").Value = 27

$ws.Range("B3").Value = 6
$ws.Range("D3").Value = 28

$ws.Range("B4").Value = 6
$ws.Range("D4").Value = 25

# Update selected cell/range on the sheet
$ws.Activate()
$ws.Range("E5").Select()
